{"js": "const replacements = [\n  [\"2024-03-23 Saturday\", \"2024-03-24 Sunday\"],\n  [\"803\u00f75=160, 3\", \"980\u00f74=245, 0\"],\n  [\"955\u00f74=238, 3\", \"808\u00f78=101, 0\"],\n  [\"715\u00f77=102, 1\", \"724\u00f75=144, 4\"],\n  [\"703\u00f75=140, 3\", \"597\u00f73=199, 0\"],\n  [\"762\u00f74=190, 2\", \"985\u00f74=246, 1\"],\n  [\"739\u00f78=92, 3\", \"274\u00f77=39, 1\"],\n  [\"662\u00f76=110, 2\", \"243\u00f77=34, 5\"],\n  [\"504\u00f79=56, 0\", \"947\u00f78=118, 3\"],\n  [\"734\u00f75=146, 4\", \"119\u00f77=17, 0\"],\n  [\"303\u00f79=33, 6\", \"452\u00f74=113, 0\"],\n  [\"287\u00f75=57, 2\", \"674\u00f76=112, 2\"],\n  [\"839\u00f74=209, 3\", \"620\u00f75=124, 0\"],\n  [\"961\u00f74=240, 1\", \"372\u00f78=46, 4\"],\n  [\"949\u00f73=316, 1\", \"864\u00f76=144, 0\"],\n  [\"975\u00f75=195, 0\", \"263\u00f79=29, 2\"],\n  [\"604\u00f73=201, 1\", \"108\u00f76=18, 0\"],\n  [\"706\u00f72=353, 0\", \"728\u00f79=80, 8\"],\n  [\"560\u00f77=80, 0\", \"225\u00f74=56, 1\"],\n  [\"778\u00f79=86, 4\", \"399\u00f76=66, 3\"],\n  [\"205\u00f78=25, 5\", \"355\u00f78=44, 3\"],\n  [\"417\u00f78=52, 1\", \"596\u00f74=149, 0\"],\n  [\"465\u00f79=51, 6\", \"713\u00f77=101, 6\"],\n  [\"287\u00f78=35, 7\", \"864\u00f78=108, 0\"],\n  [\"476\u00f75=95, 1\", \"394\u00f77=56, 2\"],\n  [\"820\u00f75=164, 0\", \"317\u00f74=79, 1\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-03-23 Saturday', '2024-03-24 Sunday'),\n    @('803\u00f75=160, 3', '980\u00f74=245, 0'),\n    @('955\u00f74=238, 3', '808\u00f78=101, 0'),\n    @('715\u00f77=102, 1', '724\u00f75=144, 4'),\n    @('703\u00f75=140, 3', '597\u00f73=199, 0'),\n    @('762\u00f74=190, 2', '985\u00f74=246, 1'),\n    @('739\u00f78=92, 3', '274\u00f77=39, 1'),\n    @('662\u00f76=110, 2', '243\u00f77=34, 5'),\n    @('504\u00f79=56, 0', '947\u00f78=118, 3'),\n    @('734\u00f75=146, 4', '119\u00f77=17, 0'),\n    @('303\u00f79=33, 6', '452\u00f74=113, 0'),\n    @('287\u00f75=57, 2', '674\u00f76=112, 2'),\n    @('839\u00f74=209, 3', '620\u00f75=124, 0'),\n    @('961\u00f74=240, 1', '372\u00f78=46, 4'),\n    @('949\u00f73=316, 1', '864\u00f76=144, 0'),\n    @('975\u00f75=195, 0', '263\u00f79=29, 2'),\n    @('604\u00f73=201, 1', '108\u00f76=18, 0'),\n    @('706\u00f72=353, 0', '728\u00f79=80, 8'),\n    @('560\u00f77=80, 0', '225\u00f74=56, 1'),\n    @('778\u00f79=86, 4', '399\u00f76=66, 3'),\n    @('205\u00f78=25, 5', '355\u00f78=44, 3'),\n    @('417\u00f78=52, 1', '596\u00f74=149, 0'),\n    @('465\u00f79=51, 6', '713\u00f77=101, 6'),\n    @('287\u00f78=35, 7', '864\u00f78=108, 0'),\n    @('476\u00f75=95, 1', '394\u00f77=56, 2'),\n    @('820\u00f75=164, 0', '317\u00f74=79, 1')\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $found = $find.Execute([ref]$pair[0], $false, $false, $false, $false, $false, $true, 1, $false, [ref]$pair[1], 2)\n    if (-not $found) {\n        throw \"No match found for: $($pair[0])\"\n    }\n}\n"}
